$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last 3 data rows (rows 53:55), shrinking the table to A1:E52
$ws.Range("A53:E55").Delete()

# New constant lambda_1 / lambda_2 values for all data rows (2..52)
$ws.Range("B2:B52").Value = 33.94444444444444
$ws.Range("C2:C52").Value = 1.95

# New dic_nbre_clients_poisson_2_keys values (column D)
$dValues = @(0,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,47,48,50,51,52,55,56,57)

# New dic_nbre_clients_prob_poisson_2_values values (column E)
$eValues = @(0.145,0.007,0.012,0.015,0.033,0.051,0.054,0.04,0.04,0.026,0.038,0.027,0.032,0.035,0.039,0.036,0.022,0.026,0.026,0.023,0.018,0.014,0.023,0.02,0.02,0.025,0.018,0.01,0.011,0.011,0.008,0.011,0.014,0.008,0.011,0.008,0.007,0.006,0.003,0.003,0.005,0.003,0.003,0.002,0.003,0.002,0.001,0.001,0.001,0.001,0.001)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
